$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells stay as Text (values look numeric e.g. "1.005", "27.530.19")
# so Excel does not auto-convert them to numbers and strip formatting / dots.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.530.19"
$ws.Range("E2").Value = "  -2.81%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.750.80"
$ws.Range("E3").Value = "  -3.52%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  +0.30%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.26"
$ws.Range("E5").Value = "  -0.83%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.003"
$ws.Range("E6").Value = "  +0.31%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4477"
$ws.Range("E7").Value = "  +2.41%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3593"
$ws.Range("E8").Value = "  -2.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07511"
$ws.Range("E9").Value = "  -2.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.28"
$ws.Range("E10").Value = "  -5.34%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.102"
$ws.Range("E11").Value = "  -3.59%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.004"
$ws.Range("E12").Value = "  +0.33%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.64"
$ws.Range("E13").Value = "  -6.18%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.189"
$ws.Range("E15").Value = "  -4.40%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.755.82"
$ws.Range("E16").Value = "  -3.55%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "92.96"
$ws.Range("E17").Value = "  -2.63%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001064"
$ws.Range("E18").Value = "  -1.59%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06410"
$ws.Range("E19").Value = "  -1.58%  "
$ws.Range("E20").Value = "  +0.23%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.03"
$ws.Range("E21").Value = "  -2.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.852"
$ws.Range("E22").Value = "  -6.40%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.578.68"
$ws.Range("E23").Value = "  -2.69%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.21"
$ws.Range("E24").Value = "  -3.21%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.115"
$ws.Range("E25").Value = "  +1.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "162.80"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.40"
$ws.Range("E27").Value = "  -1.54%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.957.93"
$ws.Range("E28").Value = "  -3.37%  "
$ws.Range("E29").Value = "  -6.64%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.78"
$ws.Range("E30").Value = "  -2.50%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.091"
$ws.Range("E31").Value = "  -9.56%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09037"
$ws.Range("E32").Value = "  -1.60%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.640"
$ws.Range("E33").Value = "  +3.64%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.538"
$ws.Range("E34").Value = "  -7.30%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "12.06"
$ws.Range("E35").Value = "  -7.19%  "
$ws.Range("E36").Value = "  -2.19%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2101"
$ws.Range("E37").Value = "  -3.35%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.6378"
$ws.Range("E38").Value = "  -3.13%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05969"
$ws.Range("E39").Value = "  -3.85%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.939"
$ws.Range("E40").Value = "  -4.96%  "
$ws.Range("E41").Value = "  -0.18%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.003"
$ws.Range("E42").Value = "  +0.30%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.391"
$ws.Range("E43").Value = "  -2.66%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.785"
$ws.Range("E44").Value = "  -4.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.23"
$ws.Range("E45").Value = "  -4.53%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.712"
$ws.Range("E46").Value = "  -0.95%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5871"
$ws.Range("E47").Value = "  -3.95%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.963"
$ws.Range("E48").Value = "  -2.65%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "121.53"
$ws.Range("E49").Value = "  -3.50%  "
$ws.Range("E50").Value = "  -0.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06865"
$ws.Range("E51").Value = "  -1.93%  "
